{"js": "// Append the admin-portal SQL reference block (select/insert statements)\n// to the end of the document body, after the existing \"400, 50 TITLE , point 22\"\n// paragraph and before the section break.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Anchor on the existing \"400, 50 TITLE , point 22\" paragraph (the last\n// paragraph in the document before this edit) and append the new content\n// right after it, before the section break.\nlet anchor = null;\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf(\"point 22\") !== -1) {\n    anchor = p;\n  }\n}\nif (!anchor) {\n  // Fallback: use the last paragraph in the body.\n  anchor = paragraphs.items[paragraphs.items.length - 1];\n}\n\nconst newLines = [\n  \"\",\n  \"select * from admin\",\n  \"select * from doctors\",\n  \"select * from Patients\",\n  \"select * from Staff\",\n  \"\",\n  \"\",\n  \"INSERT INTO Staff (branch_id, first_name, last_name, cnic, password, role, contact_number, email, address, hire_date)\",\n  \"VALUES\",\n  \"(1, 'Nadia', 'Khan', '6', 's123', 'Nurse', '03001112233', 'nadiakhan@example.com', '789 Street C', GETDATE()),\",\n  \"(1, 'Ahmed', 'Raza', '7', 's456', 'Cleaner', '03004455667', 'ahmedraza@example.com', '101 Street D', GETDATE());\",\n  \"\",\n  \"INSERT INTO Doctors (branch_id, first_name, last_name, cnic, password, specialty, contact_number, email, created_at)\",\n  \"VALUES\",\n  \"(1, 'Dr. Faisal', 'Malik', '3', 'd123', 'Cardiology', '03005566778', 'faisalmalik@example.com', GETDATE()),\",\n  \"(1, 'Dr. Hina', 'Ali', '4', 'd456', 'Neurology', '03006677889', 'hinaali@example.com', GETDATE());\",\n  \"\",\n];\n\nfor (const line of newLines) {\n  anchor = anchor.insertParagraph(line, Word.InsertLocation.after);\n}\n\nawait context.sync();\n", "ps1": "# Append the admin-portal SQL reference block (select/insert statements)\n# to the end of the document body, after the existing \"400, 50 TITLE , point 22\"\n# paragraph and before the section break.\n\n$d = $word.ActiveDocument\n\n# Confirm the anchor is the existing \"400, 50 TITLE , point 22\" paragraph\n# (the last paragraph in the document before this edit); new content is\n# appended after it, before the section break.\n$anchorRange = $d.Content\n$found = $anchorRange.Find.Execute(\"point 22\")\n\nfunction Add-Para($text) {\n    $r = $d.Paragraphs.Last.Range\n    $r.Collapse(0)\n    $r.InsertParagraphAfter()\n    if ($text -ne \"\") {\n        $r2 = $d.Paragraphs.Last.Range\n        $r2.Collapse(0)\n        $r2.InsertAfter($text)\n    }\n}\n\n$newLines = @(\n    \"\",\n    \"select * from admin\",\n    \"select * from doctors\",\n    \"select * from Patients\",\n    \"select * from Staff\",\n    \"\",\n    \"\",\n    \"INSERT INTO Staff (branch_id, first_name, last_name, cnic, password, role, contact_number, email, address, hire_date)\",\n    \"VALUES\",\n    \"(1, 'Nadia', 'Khan', '6', 's123', 'Nurse', '03001112233', 'nadiakhan@example.com', '789 Street C', GETDATE()),\",\n    \"(1, 'Ahmed', 'Raza', '7', 's456', 'Cleaner', '03004455667', 'ahmedraza@example.com', '101 Street D', GETDATE());\",\n    \"\",\n    \"INSERT INTO Doctors (branch_id, first_name, last_name, cnic, password, specialty, contact_number, email, created_at)\",\n    \"VALUES\",\n    \"(1, 'Dr. Faisal', 'Malik', '3', 'd123', 'Cardiology', '03005566778', 'faisalmalik@example.com', GETDATE()),\",\n    \"(1, 'Dr. Hina', 'Ali', '4', 'd456', 'Neurology', '03006677889', 'hinaali@example.com', GETDATE());\",\n    \"\"\n)\n\nforeach ($line in $newLines) {\n    Add-Para($line)\n}\n"}
